$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Salary for the "Sachin Thete" row (row 3) drops from 40000 to 24000
$ws.Range("D3").Value = 24000

# Sirname for the "Mona" row (row 5) changes from "Patil" to "Chavanke"
$ws.Range("C5").Value = "Chavanke"

# Leave the active selection on D7, as last left by the editor
$ws.Range("D7").Select()
